$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the cell text that previously held "SWITCHfXBR" to the new value "work_mode"
$ws.Range("B8").Value = "work_mode"

# Move/update the active cell selection to B9
$ws.Activate()
$ws.Range("B9").Select()
